$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RB")

$ws.Cells.Item(7, 1).Value = "R.Armstead"
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0

$ws.Range("J8").Select()
